$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "current" sheet: append summary row for the new date 2024-07-13
# ---------------------------------------------------------------------
$current = $wb.Worksheets.Item("current")
$current.Range("A7").Value = "'2024-07-13"
$current.Range("B7").Value = 5
$current.Range("C7").Value = 1
$current.Range("D7").Value = 3
$current.Range("E7").Value = 1

# ---------------------------------------------------------------------
# 2) "2024-07-12" sheet: tiny floating point correction on B6
# ---------------------------------------------------------------------
$sheet712 = $wb.Worksheets.Item("2024-07-12")
$sheet712.Range("B6").Value = 45485.83810637731

# ---------------------------------------------------------------------
# 3) New sheet "2024-07-13" appended at the end, mirroring the layout
#    of the other daily patient-log sheets.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2024-07-13"

# header row
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "Время"
$newSheet.Range("C1").Value = "ФИО пациента"
$newSheet.Range("D1").Value = "Врач"
$newSheet.Range("E1").Value = "Врач_Индекс"
$newSheet.Range("F1").Value = "М\Ж\Р"
$newSheet.Range("G1").Value = "Дата рождения"
$newSheet.Range("H1").Value = "Причина"
$newSheet.Range("I1").Value = "Давление"

# make sure the "Время" column uses the same datetime format as the
# sibling sheets (numFmt "yyyy-mm-dd h:mm:ss")
$newSheet.Range("B2:B6").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 2
$newSheet.Range("A2").Value = "'1"
$newSheet.Range("B2").Value = 45486.68052039352
$newSheet.Range("C2").Value = "patient test"
$newSheet.Range("D2").Value = "Karp_Kuzmin"
$newSheet.Range("E2").Value = 3
$newSheet.Range("F2").Value = "М"
$newSheet.Range("G2").Value = "'2006-07-03"
$newSheet.Range("H2").Value = "рототщто"
$newSheet.Range("I2").Value = "'54678908"

# row 3
$newSheet.Range("A3").Value = "'1"
$newSheet.Range("B3").Value = 45486.68202555556
$newSheet.Range("C3").Value = "opatient name"
$newSheet.Range("D3").Value = "Yefrem_Lebedev"
$newSheet.Range("E3").Value = 2
$newSheet.Range("F3").Value = "Ж"
$newSheet.Range("G3").Value = "'2006-07-10"
$newSheet.Range("H3").Value = "utfuygu"
$newSheet.Range("I3").Value = "'54678908"

# row 4
$newSheet.Range("A4").Value = "'2"
$newSheet.Range("B4").Value = 45486.83163686343
$newSheet.Range("C4").Value = "patient test name"
$newSheet.Range("D4").Value = "Karp_Kuzmin"
$newSheet.Range("E4").Value = 3
$newSheet.Range("F4").Value = "Ж"
$newSheet.Range("G4").Value = "'2006-07-03"
$newSheet.Range("H4").Value = "reason"
$newSheet.Range("I4").Value = "'678"

# row 5
$newSheet.Range("A5").Value = "'3"
$newSheet.Range("B5").Value = 45486.83648765046
$newSheet.Range("C5").Value = "Irina Vorontsova Klementjevna"
$newSheet.Range("D5").Value = "Karp_Kuzmin"
$newSheet.Range("E5").Value = 3
$newSheet.Range("F5").Value = "Ж"
$newSheet.Range("G5").Value = "'2002-01-08"
$newSheet.Range("H5").Value = "obshee obsledowanie"
$newSheet.Range("I5").Value = "'7890"

# row 6
$newSheet.Range("A6").Value = "'2"
$newSheet.Range("B6").Value = 45486.83918565972
$newSheet.Range("C6").Value = "Kiril Vodjanow Viktorovitch"
$newSheet.Range("D6").Value = "Yefrem_Lebedev"
$newSheet.Range("E6").Value = 2
$newSheet.Range("F6").Value = "Р"
$newSheet.Range("G6").Value = "'2024-07-01"
$newSheet.Range("H6").Value = "jalobi na kashel"
$newSheet.Range("I6").Value = "'678"

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# restore the original active sheet/tab (workbook still points at
# "2024-07-10" after the edit, per the source diff)
# ---------------------------------------------------------------------
$sheet710 = $wb.Worksheets.Item("2024-07-10")
$sheet710.Activate()
